{"js": "// Remove the pre-existing \"_GoBack\" bookmark from its old location\n// (right after the \"de-\" hyphenation, before \"skewing operation...\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Locate the end of the paragraph that currently finishes with the\n// trailing space run after \"...determine the depth of those points\n// from the camera. \" -- that's where the new commentary gets appended.\nconst anchor = context.document.body.search(\"from the camera. \", {\n  matchCase: true,\n});\nanchor.load(\"items\");\nawait context.sync();\n\nif (anchor.items.length === 0) {\n  throw new Error(\"Could not locate the target paragraph to extend.\");\n}\n\nconst paraEnd = anchor.items[0].getRange(Word.RangeLocation.end);\n\nconst newSentence =\n  \"In the coming days, we plan to actually begin implementing HMATCH, C2MODEL/LMATCH, and GMATCH algorithms in MATLAB and begin performing tests to see if they perform as expected. If not, we will continue to find other point matching algorithms that may or may not be more general-purpose than the SRI algorithms, which were designed for stereoscopic imaging\";\n\n// Insert the new sentence together with the trailing \".  \" as a single\n// run so that re-inserting the \"_GoBack\" bookmark afterwards lands in\n// the middle of a run rather than exactly on a run boundary.\nparaEnd.insertText(newSentence + \".  \", Word.InsertLocation.end);\nawait context.sync();\n\n// \".  \" (period followed by two spaces) now occurs exactly once in the\n// document, right after \"...stereoscopic imaging\". Use its start as the\n// point to re-insert the \"_GoBack\" bookmark.\nconst dot = context.document.body.search(\".  \", { matchCase: true });\ndot.load(\"items\");\nawait context.sync();\n\nconst bmPoint = dot.items[0].getRange(Word.RangeLocation.start);\nbmPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the pre-existing \"_GoBack\" bookmark from its old location\n# (right after the \"de-\" hyphenation, before \"skewing operation...\").\n$goBack = $d.Bookmarks(\"_GoBack\")\n$goBack.Delete()\n\n# Locate the end of the paragraph that currently finishes with the\n# trailing space run after \"...determine the depth of those points\n# from the camera. \" -- that's where the new commentary gets appended.\n$rng = $d.Content\n$rng.Find.Execute(\"from the camera. \")\n$rng.Collapse(0)\n\n$newSentence = \"In the coming days, we plan to actually begin implementing HMATCH, C2MODEL/LMATCH, and GMATCH algorithms in MATLAB and begin performing tests to see if they perform as expected. If not, we will continue to find other point matching algorithms that may or may not be more general-purpose than the SRI algorithms, which were designed for stereoscopic imaging\"\n\n# Insert the new sentence together with the trailing \".  \" as a single\n# run so that re-inserting the \"_GoBack\" bookmark afterwards lands in\n# the middle of a run rather than exactly on a run boundary.\n$rng.InsertAfter($newSentence + \".  \")\n\n# Re-find the freshly inserted text so we get a Range anchored back\n# into the live document (rather than reusing stale character offsets).\n$target = $d.Content\n$target.Find.Execute(\"stereoscopic imaging.\")\n$target.MoveEnd(1, -1)\n$target.Collapse(0)\n\n# Re-insert the \"_GoBack\" bookmark right after \"...stereoscopic imaging\"\n# and before the trailing \".  \".\n$d.Bookmarks.Add(\"_GoBack\", $target)\n"}
